# Update cryptocurrency Price (D) and Volume(1h) (E) columns with refreshed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.051.07'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.81%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.634.89'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.00%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '522.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.16'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("E8").Value = '  +0.55%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.654.11'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.20%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.34'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.82%  '
$ws.Range("E11").Value = '  +1.53%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.337'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.56%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.104.60'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.30%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '59.049.60'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.90%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.00'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.30%  '
$ws.Range("E17").Value = '  +0.40%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.653.63'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.51%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '347.17'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.10%  '
$ws.Range("E20").Value = '  -0.82%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.29'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.51%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.15'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.66%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '61.90'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.14%  '
$ws.Range("E25").Value = '  +0.89%  '
$ws.Range("E26").Value = '  +3.84%  '
$ws.Range("E27").Value = '  -0.40%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0805'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.27%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.13'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.51%  '
$ws.Range("E30").Value = '  -0.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.26'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.78%  '
$ws.Range("E32").Value = '  +2.80%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.97'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.77%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '150.52'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.80%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.980'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +8.99%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.00'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.14%  '
$ws.Range("E37").Value = '  +0.96%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '36.72'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.79%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.848'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.46%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.68'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.30%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.42'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.46%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '278.64'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.59%  '
$ws.Range("E43").Value = '  +0.95%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.996'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0985'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.19%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.57'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.46%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0525'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.20%  '
$ws.Range("E48").Value = '  +0.47%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.29'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.33%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.992.65'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.60%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.67'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.04%  '
